$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 21: num_customers 114 -> 115, retention_rate recalculated (115/206)
$ws.Range("C21").Value = 115
$ws.Range("E21").Value = 0.558252427184466

# Row 22: num_customers 71 -> 74, cohort_size 71 -> 74 (retention_rate stays 1)
$ws.Range("C22").Value = 74
$ws.Range("D22").Value = 74
